$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A150").Value = "IMX-USD"
$ws.Range("A151").Value = "GRT-USD"
